$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-22: Name, 2021, 2022, 2023 (blank means leave empty)
$data = @(
    @{Row=2;  A="Fiske med passiva redskap;  3a";                                        B=220; C=230; D=225},
    @{Row=3;  A="Fiske med övriga passiva redskap; 22-32";                                B=160; C=154; D=144},
    @{Row=4;  A="Fiske med kräftburar;  3a";                                              B=100; C=90;  D=97},
    @{Row=5;  A="Bottentrål havskräfta rist;  3a";                                        B=90;  C=91;  D=88},
    @{Row=6;  A="Passiva redskap (torsk); 25-32";                                         B=67;  C=65;  D=54},
    @{Row=7;  A="Bottentrål havskräfta och fisk;  3a20";                                  B=56;  C=64;  D=50},
    @{Row=8;  A="Räkfiske rist;  3a4";                                                    B=42;  C=45;  D=43},
    @{Row=9;  A="Bottentrål havskräfta och fisk;  3a21";                                  B=39;  C=45;  D=38},
    @{Row=10; A="Passiva redskap (torsk); 22-24";                                         B=33;  C=23;  D=21},
    @{Row=11; A="Fiske med finmaskig bottentrål efter pelagiska arter; 30-31";            B=29;  C=30;  D=26},
    @{Row=12; A="Räkfiske tunnel och rist;  3a4";                                         B=27;  C=27;  D=18},
    @{Row=13; A="Pelagiskt fiske med aktiva redskap (flyttrål, vad); 25-29";              B=25;  C=27;  D=29},
    @{Row=14; A="Bottentrål fisk;  3a20";                                                 B=24;  C=29;  D=32},
    @{Row=15; A="Pelagiskt fiske med aktiva redskap (flyttrål, vad); 21-24";              B=21;  C=12;  D=10},
    @{Row=16; A="Pelagiskt fiske med aktiva redskap (bottentrål);  3a204";                B=10;  C=6;   D=4},
    @{Row=17; A="Pelagiskt fiske med aktiva redskap (flyttrål, vad); 30-31";              B=9;   C=7;   D=9},
    @{Row=18; A="Pelagiskt fiske med aktiva redskap (flyttrål, vad); 3a4";                B=9;   C=9;   D=7},
    @{Row=19; A="Fiske med stormaskig bottentrål (torsk); 25-32";                         B=8;   C=6;   D=5},
    @{Row=20; A="Fiske med finmaskig bottentrål efter pelagiska arter; 25-29";            B=6;   C=7;   D=5},
    @{Row=21; A="Bottentrål fisk;  4";                                                    B=6;   C=6;   D=5},
    @{Row=22; A="Fiske med stormaskig bottentrål (torsk); 22-24";                         B=1;   C=$null; D=$null}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    if ($null -ne $item.C) { $ws.Cells.Item($r, 3).Value = $item.C }
    if ($null -ne $item.D) { $ws.Cells.Item($r, 4).Value = $item.D }
}
